# VerveStacks CHE_grids - SubRES_REZoning_Sol_Win_and_Hydro.xlsx
# The "grid_cell" helper column (AG) on the "solar" sheet got refreshed /
# re-pasted from its upstream source, which re-shuffled the order of the
# CHE_xx grid-cell labels listed against rows 4-26. Re-apply the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$gridCells = @(
    "CHE_24",
    "CHE_5",
    "CHE_8",
    "CHE_7",
    "CHE_3",
    "CHE_0",
    "CHE_20",
    "CHE_1",
    "CHE_6",
    "CHE_17",
    "CHE_19",
    "CHE_12",
    "CHE_10",
    "CHE_22",
    "CHE_11",
    "CHE_15",
    "CHE_25",
    "CHE_14",
    "CHE_18",
    "CHE_13",
    "CHE_21",
    "CHE_9",
    "CHE_4"
)

$startRow = 4
for ($i = 0; $i -lt $gridCells.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("AG$row").Value = $gridCells[$i]
}
